# Append documentation about the project's repository link to the
# end of the document.
#
# Target shape (from the diff):
#   ... so they can be able to use the app.<line break>
#   <new paragraph>
#   Below is the direct repository link to our project.<line break>
#   https://github.com/BethuelThuto/GroupAssignment.git

$d = $word.ActiveDocument
$sel = $word.Selection

# Index (1-based) of the last paragraph that already exists in the
# document, i.e. the paragraph that currently ends in "...use the app."
$lastExistingParaIndex = $d.Paragraphs.Count

# Move to the very end of the document.
$sel.EndKey(6)

# Build the new content as three fresh paragraphs first (this lets the
# Word engine stamp each one with the paragraph's full run formatting -
# sz/szCs/lang - inherited from the preceding paragraph mark). We will
# then collapse the two "marker" paragraph breaks back down into plain
# line breaks, which is exactly the structure the diff shows: a
# <w:br type="textWrapping"/> ending the first paragraph and another one
# separating the two lines of the new paragraph.
$sel.TypeParagraph()
$sel.TypeParagraph()
$sel.TypeText("Below is the direct repository link to our project.")
$sel.TypeParagraph()
$sel.TypeText("https://github.com/BethuelThuto/GroupAssignment.git")

# Turn the paragraph mark that ends the original (last pre-existing)
# paragraph into a line break, so "use the app." and the blank paragraph
# merge back into one paragraph ending in a <w:br/>.
$pA = $d.Paragraphs($lastExistingParaIndex)
$pB = $d.Paragraphs($lastExistingParaIndex + 1)
$rng1 = $d.Range($pA.Range.End - 1, $pB.Range.End)
$rng1.Find.Execute("^p", $true, $false, $false, $false, $false, $true, 1, $false, "^l", 1)

# Turn the paragraph mark that ends the "Below is..." line into a line
# break too, so it and the URL line merge into a single paragraph.
$pC = $d.Paragraphs($lastExistingParaIndex + 1)
$pD = $d.Paragraphs($lastExistingParaIndex + 2)
$rng2 = $d.Range($pC.Range.End - 1, $pD.Range.End)
$rng2.Find.Execute("^p", $true, $false, $false, $false, $false, $true, 1, $false, "^l", 1)
